$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "21-11-2025"
$ws.Range("B66").Value = "The price of gold in India today is ₹12,448 per gram for 24 karat gold, ₹11,410 per gram for 22 karat gold and ₹9,336 per gram for 18 karat gold (also called 999 gold)."

$ws.Range("A66").Borders.LineStyle = 1
$ws.Range("B66").Borders.LineStyle = 1
$ws.Range("B66").WrapText = $true
